# Update cryptos list: refresh Price (col D) and Volume(1h) (col E) values,
# and swap the Kaspa/RenderToken rows (43/44) per the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.610.79"
$ws.Range("E2").Value = "  +2.46%  "
$ws.Range("D3").Value = "1.681.04"
$ws.Range("E3").Value = "  +3.06%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "'220.53"
$ws.Range("E5").Value = "  +2.91%  "
$ws.Range("E6").Value = "  +2.48%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("D8").Value = "'30.27"
$ws.Range("E8").Value = "  +5.98%  "
$ws.Range("E9").Value = "  +2.70%  "
$ws.Range("D10").Value = "'0.0649"
$ws.Range("E10").Value = "  +6.79%  "
$ws.Range("E11").Value = "  -0.78%  "
$ws.Range("D12").Value = "1.921.33"
$ws.Range("E12").Value = "  +2.99%  "
$ws.Range("E13").Value = "  +12.57%  "
$ws.Range("E14").Value = "  +9.90%  "
$ws.Range("D15").Value = "1.687.84"
$ws.Range("E15").Value = "  +3.47%  "
$ws.Range("D16").Value = "'3.99"
$ws.Range("E16").Value = "  +3.45%  "
$ws.Range("D17").Value = "30.605.92"
$ws.Range("E17").Value = "  +2.43%  "
$ws.Range("D18").Value = "'66.54"
$ws.Range("E18").Value = "  +3.80%  "
$ws.Range("D19").Value = "'245.77"
$ws.Range("E19").Value = "  +1.46%  "
$ws.Range("D20").Value = "0.0₃0728"
$ws.Range("E20").Value = "  +4.12%  "
$ws.Range("D21").Value = "'0.999"
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("D22").Value = "'4.29"
$ws.Range("E22").Value = "  +4.35%  "
$ws.Range("D23").Value = "'10.10"
$ws.Range("E24").Value = "  +1.33%  "
$ws.Range("D25").Value = "'158.14"
$ws.Range("E25").Value = "  +0.35%  "
$ws.Range("D26").Value = "'15.96"
$ws.Range("E26").Value = "  +3.14%  "
$ws.Range("E27").Value = "  +2.69%  "
$ws.Range("D28").Value = "'6.72"
$ws.Range("E28").Value = "  +2.31%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.24%  "
$ws.Range("D30").Value = "'0.0498"
$ws.Range("E30").Value = "  +2.55%  "
$ws.Range("E31").Value = "  +2.73%  "
$ws.Range("D32").Value = "'3.50"
$ws.Range("E32").Value = "  +3.92%  "
$ws.Range("D33").Value = "1.515.32"
$ws.Range("E33").Value = "  +6.55%  "
$ws.Range("D34").Value = "'3.31"
$ws.Range("E34").Value = "  +4.82%  "
$ws.Range("E35").Value = "  +7.67%  "
$ws.Range("E36").Value = "  -0.26%  "
$ws.Range("D37").Value = "'83.69"
$ws.Range("E37").Value = "  +11.18%  "
$ws.Range("E38").Value = "  +5.81%  "
$ws.Range("D39").Value = "'0.595"
$ws.Range("E39").Value = "  +8.11%  "
$ws.Range("E40").Value = "  -3.02%  "
$ws.Range("E41").Value = "  +0.39%  "
$ws.Range("D42").Value = "'0.841"
$ws.Range("E42").Value = "  +2.02%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "'1.99"
$ws.Range("E43").Value = "  +0.21%  "
$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").Value = "'0.0500"
$ws.Range("E44").Value = "  +1.63%  "
$ws.Range("E45").Value = "  +0.38%  "
$ws.Range("D46").Value = "'0.999"
$ws.Range("E46").Value = "  -0.15%  "
$ws.Range("D47").Value = "'51.96"
$ws.Range("E47").Value = "  -1.66%  "
$ws.Range("D48").Value = "'5.56"
$ws.Range("E48").Value = "  +3.86%  "
$ws.Range("D49").Value = "1.813.98"
$ws.Range("E49").Value = "  +2.27%  "
$ws.Range("D50").Value = "'94.95"
$ws.Range("E50").Value = "  +6.70%  "
$ws.Range("E51").Value = "  +1.57%  "
